$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting name/age/hobby to B/C/D
$ws.Columns.Item(1).Insert()

# New "id" column
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 20

# hobby for row 2 changed from "tt" to "table tennis"
$ws.Range("D2").Value = "table tennis"

# Update selection to A2
$ws.Range("A2").Select()
